$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 16. This pushes the existing rows
# 16-81 down to 17-82 (row 15 stays put, and a blank row appears at 16).
$ws.Rows("16:16").Insert()

# The row that is now at 16 is blank; it must be restored to the data
# that used to live in row 15 (rows 15 and 16 were identical before the
# edit, so this reconstructs that old row).
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 44923
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 100112031
$ws.Range("G16").Value = "Poroto verde"
$ws.Range("H16").Value = "Magnum"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 27000
$ws.Range("L16").Value = 28000
$ws.Range("M16").Value = 27500
$ws.Range("N16").Value = "$/saco 25 kilos"
$ws.Range("O16").Value = "Región de O'Higgins"
$ws.Range("P16").Value = 1100
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = "Hortaliza"

# Row 15 becomes a brand new record (a new weekly price observation).
$ws.Range("D15").Value = 45063
$ws.Range("K15").Value = 33000
$ws.Range("L15").Value = 34000
$ws.Range("M15").Value = 33500
$ws.Range("N15").Value = "$/malla 25 kilos"
$ws.Range("O15").Value = "Perú"
$ws.Range("P15").Value = 1340
